$d = $word.ActiveDocument

$old1 = "It is necessary in the dependency injection scenario, but since it has only one implementation, it can have the same name as Exchange Currency. Or the name (trading) will be necessary if programmer want to demonstrate the interface is used at run time (Trading is an verb)"
$new1 = "It is not necessary as it is should be an entity mapping to data."

$r = $d.Content
$r.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
